$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1761194029850746
$ws.Range("C2").Value = 0.591044776119403
$ws.Range("J2").Value = 0.01791044776119403
$ws.Range("P2").Value = 0.1432835820895522
$ws.Range("S2").Value = 0.07164179104477612
$ws.Range("B3").Value = 0.009708737864077669
$ws.Range("C3").Value = 0.01456310679611651
$ws.Range("J3").Value = 0.03398058252427184
$ws.Range("P3").Value = 0.7961165048543689
$ws.Range("S3").Value = 0.145631067961165
$ws.Range("J4").Value = 0.09523809523809523
$ws.Range("P4").Value = 0.746031746031746
$ws.Range("S4").Value = 0.1587301587301587
$ws.Range("P5").Value = 0.6
$ws.Range("S5").Value = 0.4
$ws.Range("B6").Value = 0.07142857142857142
$ws.Range("D6").Value = 0.01984126984126984
$ws.Range("F6").Value = 0.06746031746031746
$ws.Range("J6").Value = 0.25
$ws.Range("O6").Value = 0.02777777777777778
$ws.Range("Q6").Value = 0.2182539682539683
$ws.Range("R6").Value = 0.03571428571428571
$ws.Range("S6").Value = 0.3095238095238095
$ws.Range("B7").Value = 0.08256880733944955
$ws.Range("D7").Value = 0.01834862385321101
$ws.Range("F7").Value = 0.05045871559633028
$ws.Range("J7").Value = 0.1834862385321101
$ws.Range("O7").Value = 0.01376146788990826
$ws.Range("Q7").Value = 0.2201834862385321
$ws.Range("R7").Value = 0.06880733944954129
$ws.Range("S7").Value = 0.3623853211009174
$ws.Range("B8").Value = 0.096579476861167
$ws.Range("D8").Value = 0.0261569416498994
$ws.Range("E8").Value = 0.004024144869215292
$ws.Range("F8").Value = 0.05633802816901409
$ws.Range("J8").Value = 0.096579476861167
$ws.Range("O8").Value = 0.03018108651911469
$ws.Range("Q8").Value = 0.227364185110664
$ws.Range("R8").Value = 0.08450704225352113
$ws.Range("S8").Value = 0.3782696177062374
$ws.Range("B9").Value = 0.09941520467836257
$ws.Range("D9").Value = 0.02339181286549707
$ws.Range("F9").Value = 0.04678362573099415
$ws.Range("J9").Value = 0.1286549707602339
$ws.Range("O9").Value = 0.03508771929824561
$ws.Range("Q9").Value = 0.1812865497076023
$ws.Range("R9").Value = 0.05847953216374269
$ws.Range("S9").Value = 0.4269005847953216
$ws.Range("B10").Value = 0.1138157894736842
$ws.Range("D10").Value = 0.02434210526315789
$ws.Range("E10").Value = 0.001973684210526316
$ws.Range("F10").Value = 0.07631578947368421
$ws.Range("J10").Value = 0.1131578947368421
$ws.Range("O10").Value = 0.01776315789473684
$ws.Range("Q10").Value = 0.2467105263157895
$ws.Range("R10").Value = 0.0625
$ws.Range("S10").Value = 0.343421052631579
$ws.Range("G11").Value = 0.1073446327683616
$ws.Range("J11").Value = 0.0903954802259887
$ws.Range("K11").Value = 0.1779661016949153
$ws.Range("L11").Value = 0.5988700564971752
$ws.Range("S11").Value = 0.02542372881355932
$ws.Range("G12").Value = 0.7247706422018348
$ws.Range("J12").Value = 0.2293577981651376
$ws.Range("K12").Value = 0.004587155963302753
$ws.Range("L12").Value = 0.01834862385321101
$ws.Range("S12").Value = 0.02293577981651376
$ws.Range("F13").Value = 0.02127659574468085
$ws.Range("G13").Value = 0.5957446808510638
$ws.Range("J13").Value = 0.3191489361702128
$ws.Range("S13").Value = 0.06382978723404255
$ws.Range("F15").Value = 0.01773049645390071
$ws.Range("H15").Value = 0.1702127659574468
$ws.Range("I15").Value = 0.06028368794326241
$ws.Range("J15").Value = 0.3191489361702128
$ws.Range("K15").Value = 0.07446808510638298
$ws.Range("M15").Value = 0.01773049645390071
$ws.Range("O15").Value = 0.04964539007092199
$ws.Range("S15").Value = 0.2907801418439716
$ws.Range("F16").Value = 0.02362204724409449
$ws.Range("H16").Value = 0.1496062992125984
$ws.Range("I16").Value = 0.05511811023622047
$ws.Range("J16").Value = 0.4488188976377953
$ws.Range("K16").Value = 0.1102362204724409
$ws.Range("M16").Value = 0.01574803149606299
$ws.Range("O16").Value = 0.06299212598425197
$ws.Range("S16").Value = 0.1338582677165354
$ws.Range("F17").Value = 0.02407704654895666
$ws.Range("H17").Value = 0.1797752808988764
$ws.Range("I17").Value = 0.0593900481540931
$ws.Range("J17").Value = 0.4414125200642054
$ws.Range("K17").Value = 0.09951845906902086
$ws.Range("M17").Value = 0.02407704654895666
$ws.Range("O17").Value = 0.06260032102728733
$ws.Range("S17").Value = 0.1091492776886035
$ws.Range("F18").Value = 0.02906976744186046
$ws.Range("H18").Value = 0.1686046511627907
$ws.Range("I18").Value = 0.1046511627906977
$ws.Range("J18").Value = 0.4244186046511628
$ws.Range("K18").Value = 0.09883720930232558
$ws.Range("M18").Value = 0.01162790697674419
$ws.Range("O18").Value = 0.06976744186046512
$ws.Range("S18").Value = 0.09302325581395349
$ws.Range("F19").Value = 0.01407407407407407
$ws.Range("H19").Value = 0.1992592592592592
$ws.Range("I19").Value = 0.06444444444444444
$ws.Range("J19").Value = 0.3962962962962963
$ws.Range("K19").Value = 0.117037037037037
$ws.Range("M19").Value = 0.01481481481481482
$ws.Range("N19").Value = 0.0007407407407407407
$ws.Range("O19").Value = 0.08074074074074074
$ws.Range("S19").Value = 0.1125925925925926
